# Diary update: add entries for 10/11/2024 and 11/11/2024, each followed
# by a "Continued work on Bandit Problem report." line, matching the
# author's commit (more work on the results/comparison section of the
# bandit report).

$d = $word.ActiveDocument

# Locate the last paragraph in the body ("Did more work on the bandit
# problem report.") and append a fresh paragraph after it via
# InsertParagraphAfter(), which leaves that existing paragraph completely
# untouched. Range.InsertXML() REPLACES the whole content its target
# range spans, so rather than aiming it at the original paragraph (and
# risking clobbering its text), it is aimed squarely at this brand-new,
# placeholder paragraph and swaps it out for the full block of new diary
# entries, still landing right before the closing sectPr.
$lastPara = $d.Paragraphs.Last
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newContent =
  "<w:p $wNs/>" +
  "<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>10/11/2024</w:t></w:r></w:p>" +
  "<w:p $wNs><w:r><w:t>Continued work on Bandit Problem report.</w:t></w:r></w:p>" +
  "<w:p $wNs/>" +
  "<w:p $wNs><w:r><w:t>1</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>/11/2024</w:t></w:r></w:p>" +
  "<w:p $wNs><w:r><w:t>Continued work on Bandit Problem report.</w:t></w:r></w:p>" +
  "<w:p $wNs/>"

$newRange.InsertXML($newContent) | Out-Null
